$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data added below the existing table (rows 9 and 14 of the first sheet)
$ws.Range("I9").Value = "asasd"
$ws.Range("K14").Value = "asdasd"

# Update the sheet's active selection to match the edited workbook
$ws.Range("M19").Select()
